# This edit inserts a new data row at row 248 of the worksheet, shifting all
# subsequent rows (old 248-353) down by one (new 249-354), and fills the new
# row 248 with a fresh "Pepino ensalada" price record for
# Femacal de La Calera / Región de Arica y Parinacota.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above current row 248; existing row 248 (and below)
# shifts down to row 249 (and below).
$ws.Rows.Item(248).Insert()

# Populate the newly inserted row 248 with the new record's values.
$ws.Cells.Item(248, 1).Value()  = 3
$ws.Cells.Item(248, 2).Value()  = "Femacal de La Calera"
$ws.Cells.Item(248, 3).Value()  = "Coquimbo"
$ws.Cells.Item(248, 4).Value()  = 44704
$ws.Cells.Item(248, 5).Value()  = 5
$ws.Cells.Item(248, 6).Value()  = 100112043
$ws.Cells.Item(248, 7).Value()  = "Pepino ensalada"
$ws.Cells.Item(248, 8).Value()  = "Sin especificar"
$ws.Cells.Item(248, 9).Value()  = "Primera"
$ws.Cells.Item(248, 10).Value() = 105
$ws.Cells.Item(248, 11).Value() = 18000
$ws.Cells.Item(248, 12).Value() = 18500
$ws.Cells.Item(248, 13).Value() = 18262
$ws.Cells.Item(248, 14).Value() = "$/caja 70 unidades"
$ws.Cells.Item(248, 15).Value() = "Región de Arica y Parinacota"
$ws.Cells.Item(248, 16).Value() = 261
$ws.Cells.Item(248, 17).Value() = 70
$ws.Cells.Item(248, 18).Value() = "Hortaliza"

# Copy the date cell style (s="2", date number format) from the row below
# onto the new date cell so formatting stays consistent with the rest of
# the "Fecha" column.
$ws.Cells.Item(249, 4).Copy()
$ws.Cells.Item(248, 4).PasteSpecial(-4122)
$excel.CutCopyMode = 0

Write-Host "Inserted new row 248; used range now:" $ws.UsedRange.Address()
